$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, pushing the existing data (rows 70-84) down to rows 71-85.
$ws.Rows("70:70").Insert()

# Populate the newly inserted row 70 with the new weekly record.
$ws.Cells.Item(70, 1).Value = 11
$ws.Cells.Item(70, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(70, 3).Value = "Bíobío"
$ws.Cells.Item(70, 4).Value = 44551
$ws.Cells.Item(70, 5).Value = 8
$ws.Cells.Item(70, 6).Value = 100112032
$ws.Cells.Item(70, 7).Value = "Zapallo italiano"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 100
$ws.Cells.Item(70, 11).Value = 9000
$ws.Cells.Item(70, 12).Value = 10000
$ws.Cells.Item(70, 13).Value = 9500
$ws.Cells.Item(70, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(70, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(70, 16).Value = 190
$ws.Cells.Item(70, 17).Value = 50
$ws.Cells.Item(70, 18).Value = "Hortaliza"
